$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @("09-10-2021", "10-10-2021", "11-10-2021", "12-10-2021", "13-10-2021", "14-10-2021")

$startRow = 283
$endRow = $startRow + $dates.Length - 1

# Force column A for the new rows to be treated as plain text so that
# date-like strings (e.g. "09-10-2021") are not auto-converted into date
# serial numbers by Excel's normal "smart" value parsing.
$dateRange = $ws.Range("A${startRow}:A${endRow}")
$dateRange.NumberFormat = "@"

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = 1694
    $ws.Cells.Item($row, 3).Value = 2114
    $ws.Cells.Item($row, 4).Value = 12227
    $ws.Cells.Item($row, 5).Value = 2134
    $ws.Cells.Item($row, 6).Value = 3696
    $ws.Cells.Item($row, 7).Value = 7100
}

# Restore the default (unformatted) look of the cells, matching the rest
# of the sheet, now that the text values are safely stored.
$dateRange.ClearFormats()
